# Apply the updated cryptocurrency price / 1h-volume figures described by the commit diff.
# Each text value is written with the cell pre-formatted as Text ("@") so that price
# strings which look numeric (e.g. "0.998") are NOT auto-converted to numbers by Excel,
# matching the workbook's existing inlineStr/text cells. The number format is then restored
# to General and the style reset to Normal so no stray formatting is introduced.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($cellRef, $text) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.NumberFormat = "General"
    $rng.Style = "Normal"
}

# Row 2
Set-TextCell "D2" "67.480.12"
Set-TextCell "E2" "  -1.58%  "

# Row 3
Set-TextCell "D3" "2.672.54"
Set-TextCell "E3" "  -1.01%  "

# Row 4
Set-TextCell "E4" "  -0.01%  "

# Row 5
Set-TextCell "D5" "598.45"
Set-TextCell "E5" "  -0.03%  "

# Row 6
Set-TextCell "D6" "166.19"
Set-TextCell "E6" "  +3.77%  "

# Row 7
Set-TextCell "E7" "  -0.05%  "

# Row 8
Set-TextCell "D8" "0.546"
Set-TextCell "E8" "  +0.43%  "

# Row 9
Set-TextCell "D9" "2.671.73"
Set-TextCell "E9" "  -0.98%  "

# Row 10
Set-TextCell "D10" "0.143"
Set-TextCell "E10" "  +1.58%  "

# Row 11
Set-TextCell "E11" "  +1.22%  "

# Row 12
Set-TextCell "D12" "0.358"
Set-TextCell "E12" "  -0.54%  "

# Row 13
Set-TextCell "E13" "  -1.53%  "

# Row 14
Set-TextCell "D14" "27.80"
Set-TextCell "E14" "  -1.68%  "

# Row 15
Set-TextCell "D15" "3.158.01"
Set-TextCell "E15" "  -1.04%  "

# Row 16
Set-TextCell "D16" "0.0000185"
Set-TextCell "E16" "  -1.95%  "

# Row 17
Set-TextCell "D17" "67.390.59"
Set-TextCell "E17" "  -1.82%  "

# Row 18
Set-TextCell "D18" "2.669.06"
Set-TextCell "E18" "  -0.90%  "

# Row 19
Set-TextCell "D19" "11.73"
Set-TextCell "E19" "  -1.08%  "

# Row 20
Set-TextCell "D20" "7.66"
Set-TextCell "E20" "  +0.45%  "

# Row 21
Set-TextCell "D21" "363.65"
Set-TextCell "E21" "  -0.46%  "

# Row 22
Set-TextCell "D22" "4.37"
Set-TextCell "E22" "  -3.43%  "

# Row 23
Set-TextCell "D23" "4.80"
Set-TextCell "E23" "  -1.97%  "

# Row 24
Set-TextCell "D24" "2.03"
Set-TextCell "E24" "  -4.26%  "

# Row 26
Set-TextCell "D26" "70.84"
Set-TextCell "E26" "  -4.71%  "

# Row 27
Set-TextCell "D27" "10.04"
Set-TextCell "E27" "  +1.32%  "

# Row 28
Set-TextCell "D28" "2.721.75"
Set-TextCell "E28" "  -4.11%  "

# Row 29
Set-TextCell "E29" "  -2.63%  "

# Row 30
Set-TextCell "D30" "0.998"
Set-TextCell "E30" "  -0.09%  "

# Row 31
Set-TextCell "D31" "556.62"
Set-TextCell "E31" "  -4.18%  "

# Row 32
Set-TextCell "D32" "8.00"
Set-TextCell "E32" "  -2.83%  "

# Row 33
Set-TextCell "D33" "1.39"
Set-TextCell "E33" "  -3.73%  "

# Row 34
Set-TextCell "D34" "1.93"
Set-TextCell "E34" "  -0.90%  "

# Row 35
Set-TextCell "D35" "0.129"
Set-TextCell "E35" "  -2.51%  "

# Row 37
Set-TextCell "E37" "  -5.13%  "

# Row 38
Set-TextCell "D38" "19.54"
Set-TextCell "E38" "  -1.27%  "

# Row 39
Set-TextCell "D39" "154.52"
Set-TextCell "E39" "  -4.45%  "

# Row 40
Set-TextCell "D40" "0.373"
Set-TextCell "E40" "  -1.80%  "

# Row 41
Set-TextCell "D41" "5.30"
Set-TextCell "E41" "  -1.76%  "

# Row 42
Set-TextCell "E42" "  -4.33%  "

# Row 43
Set-TextCell "D43" "17.95"
Set-TextCell "E43" "  +0.51%  "

# Row 44
Set-TextCell "E44" "  +0.01%  "

# Row 45
Set-TextCell "D45" "2.52"
Set-TextCell "E45" "  -6.11%  "

# Row 46
Set-TextCell "D46" "40.35"
Set-TextCell "E46" "  -0.73%  "

# Row 47
Set-TextCell "D47" "0.0₆0298"
Set-TextCell "E47" "  -6.05%  "

# Row 48
Set-TextCell "D48" "0.591"
Set-TextCell "E48" "  -1.67%  "

# Row 49
Set-TextCell "D49" "153.26"
Set-TextCell "E49" "  -2.90%  "

# Row 50
Set-TextCell "D50" "3.84"
Set-TextCell "E50" "  -2.77%  "

# Row 51
Set-TextCell "D51" "1.72"
Set-TextCell "E51" "  -2.85%  "

